$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1-7 (header + Logistic Regression) remain unchanged.

# New table data starting at row 8, replacing the old SVM/XGBoost rows
# and inserting a new Random Forest block before SVM, plus extra rows
# for SVM (class_weight) and XGBoost (min_child_weight).
$data = @(
    @("Random Forest", "n_estimators", "500, 1000, 1500"),
    @("Random Forest", "max_features", "sqrt"),
    @("Random Forest", "max_depth", "20, 40, 60, 80"),
    @("Random Forest", "bootstrap", "True, False"),
    @("Random Forest", "min_samples_leaf", "1, 2"),
    @("SVM", "c_values", "0.01, 0.1, 1, 10"),
    @("SVM", "kernel_grid", "rbf, poly"),
    @("SVM", "gamma_grid", "0.001, 0.01, 0.1, 1"),
    @("SVM", "degree_grid", "2, 4"),
    @("SVM", "class_weight", "balanced, None"),
    @("SVM", "k_folds", "5"),
    @("XGBoost", "n_estimators", "100, 200, 300"),
    @("XGBoost", "learning_rate", "0.1, 0.3, 0.5"),
    @("XGBoost", "max_depth", "2, 3, 5"),
    @("XGBoost", "min_child_weight", "1, 3, 5"),
    @("XGBoost", "k_folds", "5")
)

$startRow = 8
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Apply the same formatting used by the existing data rows (style index 2:
# centered, no border) to the whole new/changed range in one shot.
$ws.Range("A2:C2").Copy()
$ws.Range("A${startRow}:C${endRow}").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
